# merge #2, ProjectObject branch successfully integrated
#
# Re-sorts/re-baselines the "Resources" sheet's resource-demand table
# (rows 3-14) and recases the resource "Type" column from "Renewable"
# to "RENEWABLE". Also widens column G ("Assigned To") to fit the
# (now much longer) activity-list strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resources")

# row -> A (ID), B (Name), C (Type), D (Availability), E (Cost/Use),
#        F (Cost/Unit), G (Assigned To), H (Total Cost)
$rows = @(
    @{ Row = 3;  A = 10; B = "plasterer";           C = "RENEWABLE"; D = "4 #4"; E = 0; F = 41.38; G = "32[3 #4];";                                                                                          H = 14896.8 }
    @{ Row = 4;  A = 12; B = "joiner";               C = "RENEWABLE"; D = "5 #5"; E = 0; F = 42;    G = "40[5 #5];";                                                                                          H = 8400 }
    @{ Row = 5;  A = 4;  B = "mason";                C = "RENEWABLE"; D = "4 #4"; E = 0; F = 40.65; G = "27;45;55;";                                                                                          H = 9756 }
    @{ Row = 6;  A = 1;  B = "plumber";              C = "RENEWABLE"; D = "4 #4"; E = 0; F = 36;    G = "30;52;53;";                                                                                          H = 8640 }
    @{ Row = 7;  A = 2;  B = "electrician";          C = "RENEWABLE"; D = "4 #4"; E = 0; F = 32;    G = "29;51;61;54;";                                                                                       H = 24320 }
    @{ Row = 8;  A = 6;  B = "roofer";               C = "RENEWABLE"; D = "4 #4"; E = 0; F = 36;    G = "49[4 #4];";                                                                                          H = 11520 }
    @{ Row = 9;  A = 11; B = "painter";              C = "RENEWABLE"; D = "6 #6"; E = 0; F = 35;    G = "39[4 #6];";                                                                                          H = 16800 }
    @{ Row = 10; A = 9;  B = "screed layer";         C = "RENEWABLE"; D = "3 #3"; E = 0; F = 0;     G = "34[2 #3];56;";                                                                                       H = 0 }
    @{ Row = 11; A = 7;  B = "carpenter";            C = "RENEWABLE"; D = "3 #3"; E = 0; F = 38;    G = "26[3 #3];";                                                                                          H = 9120 }
    @{ Row = 12; A = 3;  B = "team subcontractor";   C = "RENEWABLE"; D = "8 #8"; E = 0; F = 38.56; G = "6[4 #8];17;69;18[3 #8];21[3 #8];23[5 #8];19[4 #8];20[4 #8];22[3 #8];24[4 #8];25[2 #8];28[5 #8];"; H = 173057.28 }
    @{ Row = 13; A = 5;  B = "heating engineer";     C = "RENEWABLE"; D = "2 #2"; E = 0; F = 43.2;  G = "59;60;";                                                                                             H = 13824 }
    @{ Row = 14; A = 8;  B = "floor layer";          C = "RENEWABLE"; D = "4 #4"; E = 0; F = 41.56; G = "35;57;58;";                                                                                          H = 19948.8 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
}

# Widen the "Assigned To" column (G) to fit the longer activity lists now
# held there (new explicit column width, matching a 40-character column).
$ws.Columns.Item(7).ColumnWidth = 40

# The "Activity Duration Distribution Profiles" labels on the Risk
# Analysis sheet are recased to upper snake-case to match the
# ProjectObject enum names (column D, rows 5-61).
$wsRisk = $wb.Worksheets.Item("Risk Analysis")
$profileMap = @{
    "manual - absolute"       = "MANUAL - ABSOLUTE"
    "standard - symmetric"    = "STANDARD - SYMMETRIC"
    "standard - no risk"      = "STANDARD - NO_RISK"
    "standard - skewed left"  = "STANDARD - SKEWED_LEFT"
}

for ($r = 3; $r -le 61; $r++) {
    $cell = $wsRisk.Range("D$r")
    $v = $cell.Value2
    if ($v -and $profileMap.ContainsKey($v)) {
        $cell.Value = $profileMap[$v]
    }
}
